$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refreshed prices/volumes and reordering of two row pairs
# (Stellar/WrappedEther swap rows 19-20; Kaspa/USDe swap rows 39-40)

$ws.Range("D2").Value = '96.758.73'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '3.339.24'
$ws.Range("E3").Value = '  -1.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.60'
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '655.85'
$ws.Range("E6").Value = '  +1.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.41'
$ws.Range("E7").Value = '  -2.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.423'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").Value = '  -3.83%  '
$ws.Range("D11").Value = '3.338.55'
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.49'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").Value = '96.397.35'
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.10'
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000252'
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("D17").Value = '3.959.35'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.74'
$ws.Range("E18").Value = '  +3.95%  '
$ws.Range("B19").Value = 'Stellar'
$ws.Range("C19").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.579'
$ws.Range("E19").Value = '  +20.86%  '
$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").Value = '3.320.20'
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.62'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '507.75'
$ws.Range("E23").Value = '  +1.41%  '
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000199'
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.57'
$ws.Range("E26").Value = '  +9.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '96.58'
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.11'
$ws.Range("E28").Value = '  -3.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.146'
$ws.Range("E29").Value = '  -3.51%  '
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.15'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.189'
$ws.Range("E32").Value = '  -4.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.50'
$ws.Range("E33").Value = '  +12.06%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.551'
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.31'
$ws.Range("E36").Value = '  -3.77%  '
$ws.Range("E37").Value = '  +6.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.75'
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.151'
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '509.81'
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.36'
$ws.Range("E42").Value = '  -1.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0433'
$ws.Range("E43").Value = '  +5.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.834'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.65'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("E46").Value = '  +7.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.55'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.47'
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.40'
$ws.Range("E49").Value = '  +3.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.13'
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '162.04'
$ws.Range("E51").Value = '  +1.03%  '
